# Update State GDP (column I) hardcoded values on the "OECD Data" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OECD Data")

$ws.Range("I2").Value = 309239200000
$ws.Range("I3").Value = 313086200000
$ws.Range("I4").Value = 319129900000
$ws.Range("I5").Value = 325322700000
$ws.Range("I6").Value = 333919500000
$ws.Range("I7").Value = 338752300000
